# #3220 updated BPS sample data targets
# Updates the GHGI Target (column P) and EUI Target (column R) values for the
# BPS Data sample sheet, matches the selection left active in the workbook,
# and re-syncs the conditional-formatting rule priorities that Excel
# renumbered when the rules were reviewed/reordered in the Rules Manager.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Updated GHGI Target (P) / EUI Target (R) sample values
# ---------------------------------------------------------------------------
$ws.Range("P2").Value = 60
$ws.Range("R2").Value = 100

$ws.Range("P3").Value = 30

$ws.Range("P4").Value = 45
$ws.Range("R4").Value = 68

$ws.Range("P5").Value = 30
$ws.Range("R5").Value = 56

$ws.Range("P6").Value = 45
$ws.Range("R6").Value = 68

$ws.Range("P7").Value = 30
$ws.Range("R7").Value = 56

$ws.Range("P8").Value = 40
$ws.Range("R8").Value = 60

$ws.Range("P9").Value = 40

$ws.Range("P10").Value = 38
$ws.Range("R10").Value = 44

# ---------------------------------------------------------------------------
# 2. Conditional formatting: the "K2:K10 / M2:M10" rule group (Gross Area +
#    Site EUI range checks) was moved from the top of the priority list to
#    just above the "L2:M10" group, while every other rule kept its relative
#    order. Re-assign priorities to match.
# ---------------------------------------------------------------------------
$ws.Range("L2:M10").FormatConditions.Item(1).Priority = 48   # $R2="Yes"
$ws.Range("L2:M10").FormatConditions.Item(2).Priority = 45   # $P2="Yes"

$ws.Range("P2:Q6").FormatConditions.Item(1).Priority = 5     # $T2="Yes"
$ws.Range("P2:Q6").FormatConditions.Item(2).Priority = 3     # $S2="Yes"

$ws.Range("R2").FormatConditions.Item(3).Priority = 4        # notBetween 40,1000

$ws.Range("P9:P10").FormatConditions.Item(1).Priority = 6    # $T6="Yes"
$ws.Range("P9:P10").FormatConditions.Item(2).Priority = 7    # $S6="Yes"

$ws.Range("Q10").FormatConditions.Item(3).Priority = 1       # $T8="Yes"
$ws.Range("Q10").FormatConditions.Item(4).Priority = 2       # $S8="Yes"

$ws.Range("M2").FormatConditions.Item(3).Priority = 16       # notBetween 40,1000

$ws.Range("K2:K10").FormatConditions.Item(1).Priority = 10   # $R2="Yes"
$ws.Range("K2:K10").FormatConditions.Item(2).Priority = 9    # notBetween 1000,1000000
$ws.Range("K2:K10").FormatConditions.Item(3).Priority = 8    # $P2="Yes"

# ---------------------------------------------------------------------------
# 3. Selection left on the sheet after the edits
# ---------------------------------------------------------------------------
$ws.Range("P1:S10").Select()
